$p = $ppt.ActivePresentation
Write-Output "Slide count: $($p.Slides.Count)"
